# Correction type pour génération à partir fsh ea4a6f04ed193a83290686b2f69a3f9cd2e7f4ad
#
# "Metadata" sheet (sheet1):
#  - B4 ("Name" value) was empty -> now set to "OrientationparticuliereVs"
#  - B8 ("Date" value) refreshed to the new generation timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B4").Value = "OrientationparticuliereVs"
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
